$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date and Count values ---
$metadata = $wb.Worksheets.Item("Metadata")

# "Date" value (B8) -> plain text already (not numeric-looking), so a normal
# value assignment keeps it typed as a string.
$metadata.Range("B8").Value = "2024-09-13T20:57:31+00:00"

# "Count" value (B22) looks like a number ("3"), so force text formatting
# first or Excel will store it as a numeric cell instead of a string, then
# restore the original (unformatted) look by re-applying the formatting of
# the neighboring cell that already has the right style.
$metadata.Range("B22").NumberFormat = "@"
$metadata.Range("B22").Value = "3"
$metadata.Range("B21").Copy()
$metadata.Range("B22").PasteSpecial(-4122)

# --- Concepts sheet: append a new "unknown" concept row ---
$concepts = $wb.Worksheets.Item("Concepts")

# Level (A4) is the text "1", which again looks numeric, so force text.
$concepts.Range("A4").NumberFormat = "@"
$concepts.Range("A4").Value = "1"
$concepts.Range("B4").Value = "unknown"
$concepts.Range("C4").Value = "Unknown"
$concepts.Range("D4").Value = ""

# Match the look of the row above (border/shading/alignment) for the new row.
$concepts.Range("A3:D3").Copy()
$concepts.Range("A4:D4").PasteSpecial(-4122)
